$wb = $excel.ActiveWorkbook

# ---- Sheet: Home win ----
$ws = $wb.Worksheets.Item("Home win")
$ws.Cells.Item(2, 1).Value = "26-03-2025 15:00"
$ws.Cells.Item(2, 2).Value = "SLOVAKIA"
$ws.Cells.Item(2, 3).Value = "3. LIGA - EAST"
$ws.Cells.Item(2, 4).Value = "Spišské Podhradie - Partizán Bardejov"
$ws.Cells.Item(2, 5).Value = 73.3
$ws.Cells.Item(2, 6).Value = 1.7
$ws.Cells.Item(3, 1).Value = "26-03-2025 12:00"
$ws.Cells.Item(3, 2).Value = "TURKEY"
$ws.Cells.Item(3, 3).Value = "2. LIG"
$ws.Cells.Item(3, 4).Value = "Kırklarelispor - Ankaraspor"
$ws.Cells.Item(3, 5).Value = 73.3
$ws.Cells.Item(3, 6).Value = 2.15
$ws.Cells.Item(4, 1).Value = "27-03-2025 18:45"
$ws.Cells.Item(4, 2).Value = "ROMANIA"
$ws.Cells.Item(4, 3).Value = "LIGA II"
$ws.Cells.Item(4, 4).Value = "Csikszereda - FC Voluntari"
$ws.Cells.Item(4, 5).Value = 80
$ws.Cells.Item(4, 6).Value = 2.3
$ws.Rows.Item(5).Delete()

# ---- Sheet: Btts ----
$ws = $wb.Worksheets.Item("Btts")
$ws.Cells.Item(2, 1).Value = "26-03-2025 19:00"
$ws.Cells.Item(2, 2).Value = "BRAZIL"
$ws.Cells.Item(2, 3).Value = "PAULISTA - A3"
$ws.Cells.Item(2, 4).Value = "Monte Azul - Itapirense"
$ws.Cells.Item(2, 5).Value = 76
$ws.Cells.Item(2, 6).Value = 2.2
$ws.Cells.Item(3, 1).Value = "26-03-2025 23:30"
$ws.Cells.Item(3, 2).Value = "BRAZIL"
$ws.Cells.Item(3, 3).Value = "PAULISTA - A3"
$ws.Cells.Item(3, 4).Value = "Marília - União São João"
$ws.Cells.Item(3, 5).Value = 85
$ws.Cells.Item(3, 6).Value = 2.15
$ws.Cells.Item(4, 1).Value = "26-03-2025 20:30"
$ws.Cells.Item(4, 2).Value = "EGYPT"
$ws.Cells.Item(4, 3).Value = "SECOND LEAGUE"
$ws.Cells.Item(4, 4).Value = "Kahraba Ismailia - Abu Qair Semad"
$ws.Cells.Item(4, 5).Value = 76
$ws.Cells.Item(4, 6).Value = 1.91
$ws.Cells.Item(5, 1).Value = "26-03-2025 12:00"
$ws.Cells.Item(5, 2).Value = "TURKEY"
$ws.Cells.Item(5, 3).Value = "2. LIG"
$ws.Cells.Item(5, 4).Value = "Karacabey Belediyespor - Belediye Derincespor"
$ws.Cells.Item(5, 5).Value = 84
$ws.Cells.Item(5, 6).Value = 1.91
$ws.Cells.Item(6, 1).Value = "26-03-2025 12:00"
$ws.Cells.Item(6, 2).Value = "TURKEY"
$ws.Cells.Item(6, 3).Value = "2. LIG"
$ws.Cells.Item(6, 4).Value = "Serik Belediyespor - Van BB"
$ws.Cells.Item(6, 5).Value = 81.7
$ws.Cells.Item(6, 6).Value = 1.8
$ws.Cells.Item(7, 1).Value = "27-03-2025 19:30"
$ws.Cells.Item(7, 2).Value = "ARGENTINA"
$ws.Cells.Item(7, 3).Value = "LIGA PROFESIONAL ARGENTINA"
$ws.Cells.Item(7, 4).Value = "Aldosivi - Union Santa Fe"
$ws.Cells.Item(7, 5).Value = 76.7
$ws.Cells.Item(7, 6).Value = 2.25

# ---- Sheet: Over_Under ----
$ws = $wb.Worksheets.Item("Over_Under")
$ws.Cells.Item(2, 1).Value = "26-03-2025 12:00"
$ws.Cells.Item(2, 2).Value = "TURKEY"
$ws.Cells.Item(2, 3).Value = "2. LIG"
$ws.Cells.Item(2, 4).Value = "Erzin Spor - 24 Erzincanspor"
$ws.Cells.Item(2, 5).Value = 80
$ws.Cells.Item(2, 6).Value = 1.75
$ws.Cells.Item(2, 7).Value = 60
$ws.Cells.Item(2, 8).Value = 2.8
$ws.Cells.Item(3, 1).Value = "26-03-2025 12:00"
$ws.Cells.Item(3, 2).Value = "TURKEY"
$ws.Cells.Item(3, 3).Value = "3. LIG - GROUP 4"
$ws.Cells.Item(3, 4).Value = "Kahramanmaraş İstiklalsp - Bergama Belediyespor"
$ws.Cells.Item(3, 5).Value = 80
$ws.Cells.Item(3, 6).Value = 1.6
$ws.Cells.Item(3, 7).Value = 66.7
$ws.Cells.Item(3, 8).Value = 2.6
$ws.Cells.Item(4, 1).Value = "27-03-2025 03:30"
$ws.Cells.Item(4, 2).Value = "USA"
$ws.Cells.Item(4, 3).Value = "MLS NEXT PRO"
$ws.Cells.Item(4, 4).Value = "Portland Timbers II - San Jose Earthquakes II"
$ws.Cells.Item(4, 5).Value = 85
$ws.Cells.Item(4, 6).Value = 1.77
$ws.Cells.Item(4, 7).Value = 50
$ws.Cells.Item(4, 8).Value = 3

# ---- Sheet: EV Home win ----
$ws = $wb.Worksheets.Item("EV Home win")
$ws.Cells.Item(2, 1).Value = "26-03-2025 23:00"
$ws.Cells.Item(2, 2).Value = "BRAZIL"
$ws.Cells.Item(2, 3).Value = "BRASILIENSE"
$ws.Cells.Item(2, 4).Value = "Brasiliense - Gama"
$ws.Cells.Item(2, 5).Value = 60
$ws.Cells.Item(2, 6).Value = 2
$ws.Cells.Item(2, 7).Value = 0.2
$ws.Cells.Item(3, 1).Value = "26-03-2025 19:00"
$ws.Cells.Item(3, 2).Value = "BRAZIL"
$ws.Cells.Item(3, 3).Value = "CBF BRASILEIRO U20"
$ws.Cells.Item(3, 4).Value = "Fluminense U20 - Fortaleza U20"
$ws.Cells.Item(3, 5).Value = 66.7
$ws.Cells.Item(3, 6).Value = 1.73
$ws.Cells.Item(3, 7).Value = 0.15
$ws.Cells.Item(4, 1).Value = "26-03-2025 19:00"
$ws.Cells.Item(4, 2).Value = "BRAZIL"
$ws.Cells.Item(4, 3).Value = "CBF BRASILEIRO U20"
$ws.Cells.Item(4, 4).Value = "Juventude U20 - Botafogo U20"
$ws.Cells.Item(4, 5).Value = 50
$ws.Cells.Item(4, 6).Value = 3.9
$ws.Cells.Item(4, 7).Value = 0.95
$ws.Cells.Item(5, 1).Value = "26-03-2025 19:00"
$ws.Cells.Item(5, 2).Value = "BRAZIL"
$ws.Cells.Item(5, 3).Value = "CBF BRASILEIRO U20"
$ws.Cells.Item(5, 4).Value = "Corinthians U20 - Cuiabá U20"
$ws.Cells.Item(5, 5).Value = 55.7
$ws.Cells.Item(5, 6).Value = 1.91
$ws.Cells.Item(5, 7).Value = 0.06
$ws.Cells.Item(6, 1).Value = "26-03-2025 19:00"
$ws.Cells.Item(6, 2).Value = "BRAZIL"
$ws.Cells.Item(6, 3).Value = "PAULISTA - A3"
$ws.Cells.Item(6, 4).Value = "Monte Azul - Itapirense"
$ws.Cells.Item(6, 5).Value = 50
$ws.Cells.Item(6, 6).Value = 2
$ws.Cells.Item(6, 7).Value = 0
$ws.Cells.Item(7, 1).Value = "26-03-2025 23:30"
$ws.Cells.Item(7, 2).Value = "BRAZIL"
$ws.Cells.Item(7, 3).Value = "PAULISTA - A3"
$ws.Cells.Item(7, 4).Value = "Marília - União São João"
$ws.Cells.Item(7, 5).Value = 50
$ws.Cells.Item(7, 6).Value = 1.8
$ws.Cells.Item(7, 7).Value = -0.1
$ws.Cells.Item(8, 1).Value = "26-03-2025 19:00"
$ws.Cells.Item(8, 2).Value = "BRAZIL"
$ws.Cells.Item(8, 3).Value = "PAULISTA - A3"
$ws.Cells.Item(8, 4).Value = "Desportivo Brasil - Rio Branco SP"
$ws.Cells.Item(8, 5).Value = 60
$ws.Cells.Item(8, 6).Value = 2.2
$ws.Cells.Item(8, 7).Value = 0.32
$ws.Cells.Item(9, 1).Value = "26-03-2025 22:00"
$ws.Cells.Item(9, 2).Value = "COLOMBIA"
$ws.Cells.Item(9, 3).Value = "PRIMERA A"
$ws.Cells.Item(9, 4).Value = "Alianza Petrolera - Fortaleza FC"
$ws.Cells.Item(9, 5).Value = 63.3
$ws.Cells.Item(9, 6).Value = 2.38
$ws.Cells.Item(9, 7).Value = 0.51
$ws.Cells.Item(10, 1).Value = "23-03-2025 15:30"
$ws.Cells.Item(10, 2).Value = "CROATIA"
$ws.Cells.Item(10, 3).Value = "FIRST NL"
$ws.Cells.Item(10, 4).Value = "Opatija - Rudes"
$ws.Cells.Item(10, 5).Value = 55.7
$ws.Cells.Item(10, 6).Value = 2.05
$ws.Cells.Item(10, 7).Value = 0.14
$ws.Cells.Item(11, 1).Value = "26-03-2025 21:30"
$ws.Cells.Item(11, 2).Value = "ECUADOR"
$ws.Cells.Item(11, 3).Value = "LIGA PRO SERIE B"
$ws.Cells.Item(11, 4).Value = "Guayaquil City FC - San Antonio"
$ws.Cells.Item(11, 5).Value = 60
$ws.Cells.Item(11, 6).Value = 1.85
$ws.Cells.Item(11, 7).Value = 0.11
$ws.Cells.Item(12, 1).Value = "26-03-2025 20:30"
$ws.Cells.Item(12, 2).Value = "EGYPT"
$ws.Cells.Item(12, 3).Value = "SECOND LEAGUE"
$ws.Cells.Item(12, 4).Value = "El Seka El Hadid - Proxy"
$ws.Cells.Item(12, 5).Value = 53.3
$ws.Cells.Item(12, 6).Value = 1.72
$ws.Cells.Item(12, 7).Value = -0.08
$ws.Cells.Item(13, 1).Value = "26-03-2025 19:00"
$ws.Cells.Item(13, 2).Value = "GERMANY"
$ws.Cells.Item(13, 3).Value = "REGIONALLIGA - NORD"
$ws.Cells.Item(13, 4).Value = "Phönix Lübeck - Werder Bremen II"
$ws.Cells.Item(13, 5).Value = 53.3
$ws.Cells.Item(13, 6).Value = 2
$ws.Cells.Item(13, 7).Value = 0.07
$ws.Cells.Item(14, 1).Value = "26-03-2025 12:30"
$ws.Cells.Item(14, 2).Value = "HONG-KONG"
$ws.Cells.Item(14, 3).Value = "SAPLING CUP"
$ws.Cells.Item(14, 4).Value = "Southern District - Rangers"
$ws.Cells.Item(14, 5).Value = 53.3
$ws.Cells.Item(14, 6).Value = 1.73
$ws.Cells.Item(14, 7).Value = -0.08
$ws.Cells.Item(15, 1).Value = "26-03-2025 16:00"
$ws.Cells.Item(15, 2).Value = "NIGERIA"
$ws.Cells.Item(15, 3).Value = "NPFL"
$ws.Cells.Item(15, 4).Value = "Bayelsa United - Rivers United"
$ws.Cells.Item(15, 5).Value = 50
$ws.Cells.Item(15, 6).Value = 1.95
$ws.Cells.Item(15, 7).Value = -0.03
$ws.Cells.Item(16, 1).Value = "26-03-2025 12:00"
$ws.Cells.Item(16, 2).Value = "NORWAY"
$ws.Cells.Item(16, 3).Value = "NASJONAL U19 CHAMPIONS LEAGUE"
$ws.Cells.Item(16, 4).Value = "Odd U19 - Viking U19"
$ws.Cells.Item(16, 5).Value = 60
$ws.Cells.Item(16, 6).Value = 2.62
$ws.Cells.Item(16, 7).Value = 0.57
$ws.Cells.Item(17, 1).Value = "26-03-2025 15:00"
$ws.Cells.Item(17, 2).Value = "SLOVAKIA"
$ws.Cells.Item(17, 3).Value = "3. LIGA - EAST"
$ws.Cells.Item(17, 4).Value = "Spišské Podhradie - Partizán Bardejov"
$ws.Cells.Item(17, 5).Value = 73.3
$ws.Cells.Item(17, 6).Value = 1.7
$ws.Cells.Item(17, 7).Value = 0.25
$ws.Cells.Item(18, 1).Value = "26-03-2025 15:00"
$ws.Cells.Item(18, 2).Value = "SLOVAKIA"
$ws.Cells.Item(18, 3).Value = "3. LIGA - WEST"
$ws.Cells.Item(18, 4).Value = "Spartak Myjava - Sereď"
$ws.Cells.Item(18, 5).Value = 66.7
$ws.Cells.Item(18, 6).Value = 3.8
$ws.Cells.Item(18, 7).Value = 1.53
$ws.Cells.Item(19, 1).Value = "26-03-2025 20:00"
$ws.Cells.Item(19, 2).Value = "SWITZERLAND"
$ws.Cells.Item(19, 3).Value = "1. LIGA CLASSIC - GROUP 1"
$ws.Cells.Item(19, 4).Value = "La Sarraz-Eclépens - La Chaux-de-Fonds"
$ws.Cells.Item(19, 5).Value = 51.3
$ws.Cells.Item(19, 6).Value = 1.7
$ws.Cells.Item(19, 7).Value = -0.13
$ws.Cells.Item(20, 1).Value = "26-03-2025 12:00"
$ws.Cells.Item(20, 2).Value = "TURKEY"
$ws.Cells.Item(20, 3).Value = "2. LIG"
$ws.Cells.Item(20, 4).Value = "Kırklarelispor - Ankaraspor"
$ws.Cells.Item(20, 5).Value = 73.3
$ws.Cells.Item(20, 6).Value = 2.15
$ws.Cells.Item(20, 7).Value = 0.58
$ws.Cells.Item(21, 1).Value = "26-03-2025 12:00"
$ws.Cells.Item(21, 2).Value = "TURKEY"
$ws.Cells.Item(21, 3).Value = "3. LIG - GROUP 1"
$ws.Cells.Item(21, 4).Value = "Kartal Bulvarspor - Bursaspor"
$ws.Cells.Item(21, 5).Value = 50
$ws.Cells.Item(21, 6).Value = 6
$ws.Cells.Item(21, 7).Value = 2
$ws.Cells.Item(22, 1).Value = "26-03-2025 20:30"
$ws.Cells.Item(22, 2).Value = "URUGUAY"
$ws.Cells.Item(22, 3).Value = "PRIMERA DIVISIÓN - APERTURA"
$ws.Cells.Item(22, 4).Value = "Liverpool Montevideo - Club Nacional"
$ws.Cells.Item(22, 5).Value = 53.3
$ws.Cells.Item(22, 6).Value = 5
$ws.Cells.Item(22, 7).Value = 1.66
$ws.Cells.Item(23, 1).Value = "26-03-2025 23:00"
$ws.Cells.Item(23, 2).Value = "URUGUAY"
$ws.Cells.Item(23, 3).Value = "PRIMERA DIVISIÓN - APERTURA"
$ws.Cells.Item(23, 4).Value = "Plaza Colonia - Boston River"
$ws.Cells.Item(23, 5).Value = 53.3
$ws.Cells.Item(23, 6).Value = 5.75
$ws.Cells.Item(23, 7).Value = 2.06
$ws.Cells.Item(24, 1).Value = "26-03-2025 14:00"
$ws.Cells.Item(24, 2).Value = "ZAMBIA"
$ws.Cells.Item(24, 3).Value = "SUPER LEAGUE"
$ws.Cells.Item(24, 4).Value = "Green Buffaloes - NAPSA Stars"
$ws.Cells.Item(24, 5).Value = 53.3
$ws.Cells.Item(24, 6).Value = 2.3
$ws.Cells.Item(24, 7).Value = 0.23
$ws.Cells.Item(25, 1).Value = "26-03-2025 14:00"
$ws.Cells.Item(25, 2).Value = "ZAMBIA"
$ws.Cells.Item(25, 3).Value = "SUPER LEAGUE"
$ws.Cells.Item(25, 4).Value = "Red Arrows - Nchanga Rangers"
$ws.Cells.Item(25, 5).Value = 66.7
$ws.Cells.Item(25, 6).Value = 1.7
$ws.Cells.Item(25, 7).Value = 0.13
$ws.Cells.Item(26, 1).Value = "27-03-2025 01:00"
$ws.Cells.Item(26, 2).Value = "ARGENTINA"
$ws.Cells.Item(26, 3).Value = "PRIMERA NACIONAL"
$ws.Cells.Item(26, 4).Value = "Quilmes - Colegiales"
$ws.Cells.Item(26, 5).Value = 50
$ws.Cells.Item(26, 6).Value = 1.83
$ws.Cells.Item(26, 7).Value = -0.08
$ws.Cells.Item(27, 1).Value = "27-03-2025 01:30"
$ws.Cells.Item(27, 2).Value = "BRAZIL"
$ws.Cells.Item(27, 3).Value = "COPA DO NORDESTE"
$ws.Cells.Item(27, 4).Value = "Ferroviario - Sousa"
$ws.Cells.Item(27, 5).Value = 60
$ws.Cells.Item(27, 6).Value = 1.91
$ws.Cells.Item(27, 7).Value = 0.15
$ws.Cells.Item(28, 1).Value = "27-03-2025 01:00"
$ws.Cells.Item(28, 2).Value = "BRAZIL"
$ws.Cells.Item(28, 3).Value = "RORAIMENSE"
$ws.Cells.Item(28, 4).Value = "Sao Raimundo - Monte Roraima"
$ws.Cells.Item(28, 5).Value = 53.3
$ws.Cells.Item(28, 6).Value = 3.1
$ws.Cells.Item(28, 7).Value = 0.65
$ws.Cells.Item(29, 1).Value = "27-03-2025 14:00"
$ws.Cells.Item(29, 2).Value = "CAMEROON"
$ws.Cells.Item(29, 3).Value = "ELITE ONE"
$ws.Cells.Item(29, 4).Value = "Fauve Azur Elite - Dynamo De Douala"
$ws.Cells.Item(29, 5).Value = 60
$ws.Cells.Item(29, 6).Value = 1.8
$ws.Cells.Item(29, 7).Value = 0.08
$ws.Cells.Item(30, 1).Value = "28-03-2025 00:00"
$ws.Cells.Item(30, 2).Value = "COLOMBIA"
$ws.Cells.Item(30, 3).Value = "PRIMERA B"
$ws.Cells.Item(30, 4).Value = "Huila - Quindio"
$ws.Cells.Item(30, 5).Value = 53.3
$ws.Cells.Item(30, 6).Value = 1.91
$ws.Cells.Item(30, 7).Value = 0.02
$ws.Cells.Item(31, 1).Value = "27-03-2025 02:00"
$ws.Cells.Item(31, 2).Value = "COSTA-RICA"
$ws.Cells.Item(31, 3).Value = "PRIMERA DIVISIÓN"
$ws.Cells.Item(31, 4).Value = "Puntarenas FC - Sporting San Jose"
$ws.Cells.Item(31, 5).Value = 63.3
$ws.Cells.Item(31, 6).Value = 1.83
$ws.Cells.Item(31, 7).Value = 0.16
$ws.Cells.Item(32, 1).Value = "27-03-2025 23:00"
$ws.Cells.Item(32, 2).Value = "COSTA-RICA"
$ws.Cells.Item(32, 3).Value = "PRIMERA DIVISIÓN"
$ws.Cells.Item(32, 4).Value = "Santa Ana - CS Cartagines"
$ws.Cells.Item(32, 5).Value = 60
$ws.Cells.Item(32, 6).Value = 2.62
$ws.Cells.Item(32, 7).Value = 0.57
$ws.Cells.Item(33, 1).Value = "27-03-2025 18:45"
$ws.Cells.Item(33, 2).Value = "ROMANIA"
$ws.Cells.Item(33, 3).Value = "LIGA II"
$ws.Cells.Item(33, 4).Value = "Csikszereda - FC Voluntari"
$ws.Cells.Item(33, 5).Value = 80
$ws.Cells.Item(33, 6).Value = 2.3
$ws.Cells.Item(33, 7).Value = 0.84
$ws.Cells.Item(34, 1).Value = "27-03-2025 23:00"
$ws.Cells.Item(34, 2).Value = "URUGUAY"
$ws.Cells.Item(34, 3).Value = "PRIMERA DIVISIÓN - APERTURA"
$ws.Cells.Item(34, 4).Value = "Wanderers - Progreso"
$ws.Cells.Item(34, 5).Value = 53.3
$ws.Cells.Item(34, 6).Value = 1.95
$ws.Cells.Item(34, 7).Value = 0.04

# ---- Sheet: EV Away win ----
$ws = $wb.Worksheets.Item("EV Away win")
$ws.Cells.Item(2, 1).Value = "26-03-2025 22:00"
$ws.Cells.Item(2, 2).Value = "EL-SALVADOR"
$ws.Cells.Item(2, 3).Value = "PRIMERA DIVISION"
$ws.Cells.Item(2, 4).Value = "Fuerte San Francisco - Cacahuatique"
$ws.Cells.Item(2, 5).Value = 66.7
$ws.Cells.Item(2, 6).Value = 3.1
$ws.Cells.Item(2, 7).Value = 1.07
$ws.Cells.Item(3, 1).Value = "26-03-2025 14:30"
$ws.Cells.Item(3, 2).Value = "ITALY"
$ws.Cells.Item(3, 3).Value = "SERIE D - GIRONE D"
$ws.Cells.Item(3, 4).Value = "Victor San Marino - Tuttocuoio"
$ws.Cells.Item(3, 5).Value = 60
$ws.Cells.Item(3, 6).Value = 2.95
$ws.Cells.Item(3, 7).Value = 0.77
$ws.Cells.Item(4, 1).Value = "26-03-2025 15:00"
$ws.Cells.Item(4, 2).Value = "SLOVAKIA"
$ws.Cells.Item(4, 3).Value = "3. LIGA - EAST"
$ws.Cells.Item(4, 4).Value = "Baník Kalinovo - Rimavská Sobota"
$ws.Cells.Item(4, 5).Value = 55.7
$ws.Cells.Item(4, 6).Value = 3.25
$ws.Cells.Item(4, 7).Value = 0.81
$ws.Cells.Item(5, 1).Value = "26-03-2025 15:00"
$ws.Cells.Item(5, 2).Value = "SLOVAKIA"
$ws.Cells.Item(5, 3).Value = "3. LIGA - WEST"
$ws.Cells.Item(5, 4).Value = "Inter Bratislava - Hamsik Academy"
$ws.Cells.Item(5, 5).Value = 53.3
$ws.Cells.Item(5, 6).Value = 5.25
$ws.Cells.Item(5, 7).Value = 1.8
$ws.Cells.Item(6, 1).Value = "26-03-2025 12:00"
$ws.Cells.Item(6, 2).Value = "TURKEY"
$ws.Cells.Item(6, 3).Value = "2. LIG"
$ws.Cells.Item(6, 4).Value = "Somaspor - Menemen Belediyespor"
$ws.Cells.Item(6, 5).Value = 56.7
$ws.Cells.Item(6, 6).Value = 1.85
$ws.Cells.Item(6, 7).Value = 0.05
$ws.Cells.Item(7, 1).Value = "27-03-2025 21:30"
$ws.Cells.Item(7, 2).Value = "PERU"
$ws.Cells.Item(7, 3).Value = "PRIMERA DIVISIÓN"
$ws.Cells.Item(7, 4).Value = "Ayacucho FC - FBC Melgar"
$ws.Cells.Item(7, 5).Value = 53.3
$ws.Cells.Item(7, 6).Value = 1.7
$ws.Cells.Item(7, 7).Value = -0.09
$ws.Cells.Item(8, 1).Value = "27-03-2025 03:30"
$ws.Cells.Item(8, 2).Value = "USA"
$ws.Cells.Item(8, 3).Value = "MLS NEXT PRO"
$ws.Cells.Item(8, 4).Value = "Portland Timbers II - San Jose Earthquakes II"
$ws.Cells.Item(8, 5).Value = 60
$ws.Cells.Item(8, 6).Value = 2.2
$ws.Cells.Item(8, 7).Value = 0.32
$ws.Range("A9:A15").EntireRow.Delete()

# ---- Sheet: EV Over 2.5 ----
$ws = $wb.Worksheets.Item("EV Over 2.5")
$ws.Cells.Item(2, 1).Value = "26-03-2025 19:00"
$ws.Cells.Item(2, 2).Value = "GERMANY"
$ws.Cells.Item(2, 3).Value = "REGIONALLIGA - NORDOST"
$ws.Cells.Item(2, 4).Value = "FSV Zwickau - FC Rot-Weiß Erfurt"
$ws.Cells.Item(2, 5).Value = 71.8
$ws.Cells.Item(2, 6).Value = 1.8
$ws.Cells.Item(2, 7).Value = 0.29
$ws.Cells.Item(3, 1).Value = "26-03-2025 12:00"
$ws.Cells.Item(3, 2).Value = "TURKEY"
$ws.Cells.Item(3, 3).Value = "2. LIG"
$ws.Cells.Item(3, 4).Value = "Erzin Spor - 24 Erzincanspor"
$ws.Cells.Item(3, 5).Value = 80
$ws.Cells.Item(3, 6).Value = 1.75
$ws.Cells.Item(3, 7).Value = 0.4
$ws.Cells.Item(4, 1).Value = "27-03-2025 21:00"
$ws.Cells.Item(4, 2).Value = "ENGLAND"
$ws.Cells.Item(4, 3).Value = "LEAGUE ONE"
$ws.Cells.Item(4, 4).Value = "Leyton Orient - Stevenage"
$ws.Cells.Item(4, 5).Value = 70
$ws.Cells.Item(4, 6).Value = 2.35
$ws.Cells.Item(4, 7).Value = 0.64
$ws.Cells.Item(5, 1).Value = "27-03-2025 03:30"
$ws.Cells.Item(5, 2).Value = "USA"
$ws.Cells.Item(5, 3).Value = "MLS NEXT PRO"
$ws.Cells.Item(5, 4).Value = "Portland Timbers II - San Jose Earthquakes II"
$ws.Cells.Item(5, 5).Value = 85
$ws.Cells.Item(5, 6).Value = 1.77
$ws.Cells.Item(5, 7).Value = 0.5
$ws.Range("A6:A7").EntireRow.Delete()

# ---- Sheet: EV Btts ----
$ws = $wb.Worksheets.Item("EV Btts")
$ws.Cells.Item(2, 1).Value = "26-03-2025 19:00"
$ws.Cells.Item(2, 2).Value = "BRAZIL"
$ws.Cells.Item(2, 3).Value = "PAULISTA - A3"
$ws.Cells.Item(2, 4).Value = "Monte Azul - Itapirense"
$ws.Cells.Item(2, 5).Value = 76
$ws.Cells.Item(2, 6).Value = 2.2
$ws.Cells.Item(2, 7).Value = 0.67
$ws.Cells.Item(3, 1).Value = "26-03-2025 23:30"
$ws.Cells.Item(3, 2).Value = "BRAZIL"
$ws.Cells.Item(3, 3).Value = "PAULISTA - A3"
$ws.Cells.Item(3, 4).Value = "Marília - União São João"
$ws.Cells.Item(3, 5).Value = 85
$ws.Cells.Item(3, 6).Value = 2.15
$ws.Cells.Item(3, 7).Value = 0.83
$ws.Cells.Item(4, 1).Value = "27-03-2025 00:30"
$ws.Cells.Item(4, 2).Value = "COLOMBIA"
$ws.Cells.Item(4, 3).Value = "PRIMERA A"
$ws.Cells.Item(4, 4).Value = "Junior - Union Magdalena"
$ws.Cells.Item(4, 5).Value = 73.3
$ws.Cells.Item(4, 6).Value = 2
$ws.Cells.Item(4, 7).Value = 0.47
$ws.Cells.Item(5, 1).Value = "26-03-2025 20:30"
$ws.Cells.Item(5, 2).Value = "EGYPT"
$ws.Cells.Item(5, 3).Value = "SECOND LEAGUE"
$ws.Cells.Item(5, 4).Value = "Kahraba Ismailia - Abu Qair Semad"
$ws.Cells.Item(5, 5).Value = 76
$ws.Cells.Item(5, 6).Value = 1.91
$ws.Cells.Item(5, 7).Value = 0.45
$ws.Cells.Item(6, 1).Value = "26-03-2025 12:00"
$ws.Cells.Item(6, 2).Value = "TURKEY"
$ws.Cells.Item(6, 3).Value = "2. LIG"
$ws.Cells.Item(6, 4).Value = "Erzin Spor - 24 Erzincanspor"
$ws.Cells.Item(6, 5).Value = 70
$ws.Cells.Item(6, 6).Value = 2
$ws.Cells.Item(6, 7).Value = 0.4
$ws.Cells.Item(7, 1).Value = "26-03-2025 12:00"
$ws.Cells.Item(7, 2).Value = "TURKEY"
$ws.Cells.Item(7, 3).Value = "2. LIG"
$ws.Cells.Item(7, 4).Value = "Buca FK - Nazilli Belediyespor"
$ws.Cells.Item(7, 5).Value = 66.7
$ws.Cells.Item(7, 6).Value = 2.25
$ws.Cells.Item(7, 7).Value = 0.5
$ws.Cells.Item(8, 1).Value = "26-03-2025 12:00"
$ws.Cells.Item(8, 2).Value = "TURKEY"
$ws.Cells.Item(8, 3).Value = "2. LIG"
$ws.Cells.Item(8, 4).Value = "Karacabey Belediyespor - Belediye Derincespor"
$ws.Cells.Item(8, 5).Value = 84
$ws.Cells.Item(8, 6).Value = 1.91
$ws.Cells.Item(8, 7).Value = 0.6
$ws.Cells.Item(9, 1).Value = "26-03-2025 12:00"
$ws.Cells.Item(9, 2).Value = "TURKEY"
$ws.Cells.Item(9, 3).Value = "2. LIG"
$ws.Cells.Item(9, 4).Value = "Serik Belediyespor - Van BB"
$ws.Cells.Item(9, 5).Value = 81.7
$ws.Cells.Item(9, 6).Value = 1.8
$ws.Cells.Item(9, 7).Value = 0.47
$ws.Cells.Item(10, 1).Value = "27-03-2025 19:30"
$ws.Cells.Item(10, 2).Value = "ARGENTINA"
$ws.Cells.Item(10, 3).Value = "LIGA PROFESIONAL ARGENTINA"
$ws.Cells.Item(10, 4).Value = "Aldosivi - Union Santa Fe"
$ws.Cells.Item(10, 5).Value = 76.7
$ws.Cells.Item(10, 6).Value = 2.25
$ws.Cells.Item(10, 7).Value = 0.73
$ws.Cells.Item(11, 1).Value = "27-03-2025 01:00"
$ws.Cells.Item(11, 2).Value = "ARGENTINA"
$ws.Cells.Item(11, 3).Value = "TORNEO FEDERAL A"
$ws.Cells.Item(11, 4).Value = "Sportivo Las Parejas - 9 De Julio Rafaela"
$ws.Cells.Item(11, 5).Value = 70
$ws.Cells.Item(11, 6).Value = 2.25
$ws.Cells.Item(11, 7).Value = 0.57
$ws.Cells.Item(12, 1).Value = "28-03-2025 00:00"
$ws.Cells.Item(12, 2).Value = "COLOMBIA"
$ws.Cells.Item(12, 3).Value = "PRIMERA A"
$ws.Cells.Item(12, 4).Value = "Once Caldas - Llaneros"
$ws.Cells.Item(12, 5).Value = 68
$ws.Cells.Item(12, 6).Value = 2.1
$ws.Cells.Item(12, 7).Value = 0.43
$ws.Cells.Item(13, 1).Value = "27-03-2025 23:00"
$ws.Cells.Item(13, 2).Value = "COSTA-RICA"
$ws.Cells.Item(13, 3).Value = "PRIMERA DIVISIÓN"
$ws.Cells.Item(13, 4).Value = "Santa Ana - CS Cartagines"
$ws.Cells.Item(13, 5).Value = 73.3
$ws.Cells.Item(13, 6).Value = 1.91
$ws.Cells.Item(13, 7).Value = 0.4

